$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '68.076.75'
$ws.Range('E2').Value = '  -0.16%  '

$ws.Range('D3').Value = '3.609.73'
$ws.Range('E3').Value = '  -1.44%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '585.64'
$ws.Range('E5').Value = '  -1.54%  '

$ws.Range('D6').Value = '193.01'
$ws.Range('E6').Value = '  +0.77%  '

$ws.Range('D7').Value = '3.602.32'
$ws.Range('E7').Value = '  -1.53%  '

$ws.Range('D8').Value = '0.619'
$ws.Range('E8').Value = '  -0.22%  '

$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('D10').Value = '0.678'
$ws.Range('E10').Value = '  -2.76%  '

$ws.Range('D11').Value = '0.151'
$ws.Range('E11').Value = '  -0.81%  '

$ws.Range('D12').Value = '55.28'
$ws.Range('E12').Value = '  -2.84%  '

$ws.Range('D13').Value = '0.0000288'
$ws.Range('E13').Value = '  +6.21%  '

$ws.Range('D14').Value = '9.99'
$ws.Range('E14').Value = '  -2.60%  '

$ws.Range('D15').Value = '4.190.32'
$ws.Range('E15').Value = '  -1.43%  '

$ws.Range('D16').Value = '3.620.65'
$ws.Range('E16').Value = '  -1.29%  '

$ws.Range('E17').Value = '  -0.39%  '

$ws.Range('D18').Value = '12.51'
$ws.Range('E18').Value = '  -0.66%  '

$ws.Range('D19').Value = '67.975.43'
$ws.Range('E19').Value = '  -0.09%  '

$ws.Range('D20').Value = '18.48'
$ws.Range('E20').Value = '  -2.01%  '

$ws.Range('E21').Value = '  -2.28%  '

$ws.Range('D22').Value = '404.60'
$ws.Range('E22').Value = '  +0.00%  '

$ws.Range('D23').Value = '13.38'
$ws.Range('E23').Value = '  +22.51%  '

$ws.Range('D24').Value = '4.25'
$ws.Range('E24').Value = '  -4.20%  '

$ws.Range('D25').Value = '85.85'
$ws.Range('E25').Value = '  -2.83%  '

$ws.Range('D26').Value = '2.95'
$ws.Range('E26').Value = '  +0.23%  '

$ws.Range('D27').Value = '3.93'
$ws.Range('E27').Value = '  +7.26%  '

$ws.Range('D28').Value = '12.57'
$ws.Range('E28').Value = '  +0.44%  '

$ws.Range('D29').Value = '6.12'
$ws.Range('E29').Value = '  +0.63%  '

$ws.Range('D30').Value = '8.14'
$ws.Range('E30').Value = '  +13.89%  '

$ws.Range('D31').Value = '9.14'
$ws.Range('E31').Value = '  -1.85%  '

$ws.Range('D32').Value = '31.51'
$ws.Range('E32').Value = '  -1.09%  '

$ws.Range('D33').Value = '676.74'
$ws.Range('E33').Value = '  +11.65%  '

$ws.Range('D34').Value = '12.22'
$ws.Range('E34').Value = '  -0.21%  '

$ws.Range('D35').Value = '0.117'
$ws.Range('E35').Value = '  +1.21%  '

$ws.Range('D36').Value = '64.46'
$ws.Range('E36').Value = '  -5.22%  '

$ws.Range('D37').Value = '42.55'
$ws.Range('E37').Value = '  -3.53%  '

$ws.Range('E38').Value = '  +8.09%  '

$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.27%  '

$ws.Range('D40').Value = '0.0₃0785'
$ws.Range('E40').Value = '  +1.96%  '

$ws.Range('D41').Value = '2.95'
$ws.Range('E41').Value = '  +17.61%  '

$ws.Range('D42').Value = '3.206.80'
$ws.Range('E42').Value = '  +15.52%  '

$ws.Range('D43').Value = '3.12'
$ws.Range('E43').Value = '  +7.93%  '

$ws.Range('E44').Value = '  -0.55%  '

$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  -0.22%  '

$ws.Range('D46').Value = '0.0421'
$ws.Range('E46').Value = '  -0.76%  '

$ws.Range('E47').Value = '  -2.68%  '

$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '3.13'
$ws.Range('E48').Value = '  -3.24%  '

$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = '8.78'
$ws.Range('E49').Value = '  -1.20%  '

$ws.Range('D50').Value = '142.94'
$ws.Range('E50').Value = '  -0.53%  '

$ws.Range('D51').Value = '2.55'
$ws.Range('E51').Value = '  -0.99%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
